$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-cluster-name string: "Resolving-Mac" -> "Neutrophils" and
# recompute every data row (2-10) with the refreshed TPM-derived NATMI values.
# Rows 11-13 (the old "Inflammatory-Mac" target-cluster rows) are removed below.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pspn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4815393333333333
$ws.Range("H2").Value = 1.444618
$ws.Range("I2").Value = 0.3617347224948818
$ws.Range("J2").Value = 0.3617347224948818
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1389376666666667
$ws.Range("N2").Value = 0.416813
$ws.Range("O2").Value = 0.01722256533596611
$ws.Range("P2").Value = 0.01722256533596611
$ws.Range("Q2").Value = 0.06690395138155555
$ws.Range("R2").Value = 0.602135562434
$ws.Range("S2").Value = 0.00622999989245567
$ws.Range("T2").Value = 0.006229999892455672

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pspn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.4815393333333333
$ws.Range("H3").Value = 1.444618
$ws.Range("I3").Value = 0.3617347224948818
$ws.Range("J3").Value = 0.3617347224948818
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.220039333333333
$ws.Range("N3").Value = 12.660118
$ws.Range("O3").Value = 0.5231115858095611
$ws.Range("P3").Value = 0.5231115858095611
$ws.Range("Q3").Value = 2.032114927213778
$ws.Range("R3").Value = 18.289034344924
$ws.Range("S3").Value = 0.1892276243266791
$ws.Range("T3").Value = 0.1892276243266791

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pspn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.4815393333333333
$ws.Range("H4").Value = 1.444618
$ws.Range("I4").Value = 0.3617347224948818
$ws.Range("J4").Value = 0.3617347224948818
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.708210666666667
$ws.Range("N4").Value = 11.124632
$ws.Range("O4").Value = 0.4596658488544727
$ws.Range("P4").Value = 0.4596658488544728
$ws.Range("Q4").Value = 1.785649292286222
$ws.Range("R4").Value = 16.070843630576
$ws.Range("S4").Value = 0.1662770982757469
$ws.Range("T4").Value = 0.166277098275747

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pspn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.259826
$ws.Range("H5").Value = 0.779478
$ws.Range("I5").Value = 0.1951825728468463
$ws.Range("J5").Value = 0.1951825728468463
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1389376666666667
$ws.Range("N5").Value = 0.416813
$ws.Range("O5").Value = 0.01722256533596611
$ws.Range("P5").Value = 0.01722256533596611
$ws.Range("Q5").Value = 0.03609961817933333
$ws.Range("R5").Value = 0.324896563614
$ws.Range("S5").Value = 0.003361544613296776
$ws.Range("T5").Value = 0.003361544613296776

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pspn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.259826
$ws.Range("H6").Value = 0.779478
$ws.Range("I6").Value = 0.1951825728468463
$ws.Range("J6").Value = 0.1951825728468463
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.220039333333333
$ws.Range("N6").Value = 12.660118
$ws.Range("O6").Value = 0.5231115858095611
$ws.Range("P6").Value = 0.5231115858095611
$ws.Range("Q6").Value = 1.096475939822667
$ws.Range("R6").Value = 9.868283458404
$ws.Range("S6").Value = 0.102102265204304
$ws.Range("T6").Value = 0.102102265204304

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pspn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.259826
$ws.Range("H7").Value = 0.779478
$ws.Range("I7").Value = 0.1951825728468463
$ws.Range("J7").Value = 0.1951825728468463
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.708210666666667
$ws.Range("N7").Value = 11.124632
$ws.Range("O7").Value = 0.4596658488544727
$ws.Range("P7").Value = 0.4596658488544728
$ws.Range("Q7").Value = 0.9634895446773333
$ws.Range("R7").Value = 8.671405902096
$ws.Range("S7").Value = 0.08971876302924559
$ws.Range("T7").Value = 0.08971876302924559

# Row 8
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Pspn"
$ws.Range("C8").Value = "Gfra1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.5898293333333333
$ws.Range("H8").Value = 1.769488
$ws.Range("I8").Value = 0.4430827046582718
$ws.Range("J8").Value = 0.4430827046582719
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1389376666666667
$ws.Range("N8").Value = 0.416813
$ws.Range("O8").Value = 0.01722256533596611
$ws.Range("P8").Value = 0.01722256533596611
$ws.Range("Q8").Value = 0.08194951130488888
$ws.Range("R8").Value = 0.737545601744
$ws.Range("S8").Value = 0.007631020830213661
$ws.Range("T8").Value = 0.007631020830213663

# Row 9
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Pspn"
$ws.Range("C9").Value = "Gfra1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.5898293333333333
$ws.Range("H9").Value = 1.769488
$ws.Range("I9").Value = 0.4430827046582718
$ws.Range("J9").Value = 0.4430827046582719
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.220039333333333
$ws.Range("N9").Value = 12.660118
$ws.Range("O9").Value = 0.5231115858095611
$ws.Range("P9").Value = 0.5231115858095611
$ws.Range("Q9").Value = 2.489102986620444
$ws.Range("R9").Value = 22.401926879584
$ws.Range("S9").Value = 0.231781696278578
$ws.Range("T9").Value = 0.231781696278578

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Pspn"
$ws.Range("C10").Value = "Gfra1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.5898293333333333
$ws.Range("H10").Value = 1.769488
$ws.Range("I10").Value = 0.4430827046582718
$ws.Range("J10").Value = 0.4430827046582719
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.708210666666667
$ws.Range("N10").Value = 11.124632
$ws.Range("O10").Value = 0.4596658488544727
$ws.Range("P10").Value = 0.4596658488544728
$ws.Range("Q10").Value = 2.187211425379556
$ws.Range("R10").Value = 19.684902828416
$ws.Range("S10").Value = 0.2036699875494802
$ws.Range("T10").Value = 0.2036699875494802

# Remove the now-obsolete rows 11-13 (previously the "Inflammatory-Mac" target-cluster rows)
$ws.Range("A11:T13").Delete()
